$wb = $excel.ActiveWorkbook
$win = $excel.ActiveWindow
$win.Height = 12576
